$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Relabel the two header cells in row 1 ("localization work")
#   B1: "green piece (cm)" -> "pcb (cm)"
#   D1: "wood  (cm)"       -> "wood "
$ws.Range("B1").Value = "pcb (cm)"
$ws.Range("D1").Value = "wood "

# Update the selection: active cell D3, selected range D3:D23
$ws.Range("D3:D23").Select()
